$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 94; this shifts the existing row 94..211
# down to 95..212 and keeps everything else intact.
$ws.Rows.Item(94).Insert()

# Populate the newly inserted row 94 with the new weekly price record.
$ws.Range("A94").Value = 10
$ws.Range("B94").Value = "Vega Modelo de Temuco"
$ws.Range("C94").Value = "La Araucanía"
$ws.Range("D94").Value = 45174
$ws.Range("D94").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E94").Value = 9
$ws.Range("F94").Value = 100112035
$ws.Range("G94").Value = "Bruselas (repollito)"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 40
$ws.Range("K94").Value = 28000
$ws.Range("L94").Value = 28000
$ws.Range("M94").Value = 28000
$ws.Range("N94").Value = "$/malla 15 kilos"
$ws.Range("O94").Value = "Región Metropolitana"
$ws.Range("P94").Value = 1867
$ws.Range("Q94").Value = 15
$ws.Range("R94").Value = "Hortaliza"
